$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$teamIds = @(
    "0eeb011c-24fb-4476-91f7-d8e28ae49c2f",
    "5e1a20f6-82bf-4dee-aa79-41702d9feb41",
    "6afc31f0-3916-443a-92c4-b5eb425a9bc3",
    "61f0969e-22a4-4374-8588-d6511915b05e",
    "7fc75193-58a0-4e7d-ab42-382ec10a8be4",
    "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7",
    "e5ca6e2b-5f54-4acd-ad7b-03e631313986",
    "0435a227-38e8-494e-b1bf-271b00893eae",
    "e5c40d19-b03a-4f5a-82c8-25540cd45e07",
    "e0228b4f-7807-45db-a3f6-8c6e1f4adf41",
    "630f61e8-543f-46e2-af63-2b62e8bc4fd2",
    "94742748-e7ab-454b-8ff6-9893440bd059",
    "57ebee16-96d2-46a6-ab16-2476b305fd91",
    "3b8adc57-0f6f-482c-8306-9830e819d666",
    "e42288a3-b5af-4464-bc45-85d438bcea11",
    "4279fd55-c2c1-440d-abaa-430f3c27be44",
    "5184566d-523a-4432-848d-ac234ffb6ac6",
    "de503c24-f17d-47a9-9a47-6f0a194f8c9c",
    "850a92da-c3d6-4fb9-a510-99626e9ad312",
    "718c6b8f-7c00-4bcb-b53c-8f3f42154362",
    "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd",
    "e7bb31c1-e095-453b-95ff-565ea62efb0a",
    "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
)

for ($i = 0; $i -lt $teamIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $teamIds[$i]
}
